$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" for the first report row.
# Shared between the Overview sheet and the de-de sheet's
# "Correspond Handoff Datetime" column, so both must be updated.
$wsOverview.Range("G2").Value = "2017-02-22 07:57:20"
$wsDeDe.Range("H2").Value = "2017-02-22 07:57:20"

# zh-cn: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2017-02-22 07:57:01"
$wsZhCn.Range("L2").Value = "2017-02-22 07:57:56"

# de-de: Correspond Handback DateTime
$wsDeDe.Range("L2").Value = "2017-02-22 07:58:19"
